# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 159 (Vega Monumental
# Concepción - Mango), pushing the existing rows 159:187 down to 160:188.

$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

$ws.Rows("159:159").Insert()

$ws.Range("A159").Value = 11
$ws.Range("B159").Value = "Vega Monumental Concepción"
$ws.Range("C159").Value = "Bíobío"
$ws.Range("D159").Value = 45154
$ws.Range("E159").Value = 8
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100108
$ws.Range("H159").Value = "Tropicales y subtropicales"
$ws.Range("I159").Value = 100108002
$ws.Range("J159").Value = "Mango"
$ws.Range("K159").Value = "Sin especificar"
$ws.Range("L159").Value = "Primera"
$ws.Range("M159").Value = 200
$ws.Range("N159").Value = 8000
$ws.Range("O159").Value = 8500
$ws.Range("P159").Value = 8250
$ws.Range("Q159").Value = "$/bandeja 4 kilos"
$ws.Range("R159").Value = "Brasil"
$ws.Range("S159").Value = 2062
$ws.Range("T159").Value = 4
